# Append the 2025-01-17 09:00 resale-number snapshot as row 30.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text-like columns -------------------------------------------------
# A30 ("2025-01-17") and D30 ("02") look like a date / a number to the
# smart-typing engine, so force them through the Text number format while
# the value is assigned, then restore the cell's default style so no
# explicit style index is left behind on the new cells (matches the rest
# of the data rows, which carry no "s" attribute).
$ws.Range("A30").NumberFormat = "@"
$ws.Range("A30").Value = "2025-01-17"
$ws.Range("A30").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "02"
$ws.Range("D30").Style = "Normal"

# Plain text columns — no ambiguity, smart-typing leaves them as text.
$ws.Range("B30").Value = "09:00:22"
$ws.Range("C30").Value = "Friday"

# --- Numeric columns -----------------------------------------------------
$ws.Range("E30").Value = 126764
$ws.Range("F30").Value = 141738
$ws.Range("G30").Value = 169140
$ws.Range("H30").Value = 158122
$ws.Range("I30").Value = -1
$ws.Range("J30").Value = 142901
$ws.Range("K30").Value = -1
$ws.Range("L30").Value = -1
$ws.Range("M30").Value = 192119
$ws.Range("N30").Value = 115513
$ws.Range("O30").Value = 45315
$ws.Range("P30").Value = 28480
$ws.Range("Q30").Value = 65530
$ws.Range("R30").Value = -1
$ws.Range("S30").Value = 48806
$ws.Range("T30").Value = -1
